# This script applies updated cryptocurrency price/volume figures
# to columns D (Price) and E (Volume 1h) as part of the scheduled
# "Updated cryptos list" GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep a Text format so that values such as
# "1.000" or "0.9999" are preserved exactly as strings rather than
# being re-interpreted as numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.244.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.855.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4637"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3714"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07292"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8873"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.12"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07867"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.810.90"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.398"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.523"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.09"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008936"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.260.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.085"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.53"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.140.17"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.98%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.42"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.45"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.051"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.01"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.061"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08810"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.140"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7692"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.168"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.514"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.721"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.105"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.57%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05222"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.936"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.045"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5125"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1629"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.439"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4798"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.08"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.645"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06207"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.15%  "
